$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the updated crypto market data (price + 1h volume change) scraped
# on Thu Mar 16 04:33:31 UTC 2023. All target cells are plain text in the
# source workbook (inline strings), so force text format before assignment
# to avoid Excel auto-converting numeric-looking strings (e.g. "1.000",
# "24.396.02") into floating point numbers / dates.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '24.396.02'
$ws.Range('E2').NumberFormat = "@"
$ws.Range('E2').Value = '  -1.80%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.651.45'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  -3.37%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.004'
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '311.12'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  -0.19%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '1.001'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +0.20%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.3639'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  -3.30%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '46.73'
$ws.Range('E8').NumberFormat = "@"
$ws.Range('E8').Value = '  -5.93%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.3244'
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  -6.10%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  -7.58%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07009'
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  -7.32%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '1.001'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +0.07%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '5.937'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  -6.10%  '
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  -8.62%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '6.573'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  -7.09%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '1.645.30'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  -3.87%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.00001034'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  -8.99%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '0.06614'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  -1.65%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '1.000'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.12%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '78.22'
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -8.02%  '
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -7.66%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '15.56'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  -10.50%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '12.48'
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '24.375.59'
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  -1.86%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '2.478'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +0.81%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.306'
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  -17.91%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '146.69'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  -3.65%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '18.50'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  -10.02%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '1.829.67'
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  -3.71%  '
$ws.Range('B30').NumberFormat = "@"
$ws.Range('B30').Value = 'BitcoinCash'
$ws.Range('C30').NumberFormat = "@"
$ws.Range('C30').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '123.71'
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  -6.96%  '
$ws.Range('B31').NumberFormat = "@"
$ws.Range('B31').Value = 'ImmutableX'
$ws.Range('C31').NumberFormat = "@"
$ws.Range('C31').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.181'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  -5.73%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.078'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  -3.40%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.621'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -19.39%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.08437'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  -4.51%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.661'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  -7.98%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '12.13'
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -12.81%  '
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  -8.56%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '1.243'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  -3.17%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.05990'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  -10.57%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.02213'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  -8.55%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.2055'
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  -8.39%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '8.091'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  -14.01%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.001'
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +0.17%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.5868'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  -9.29%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '3.758'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -1.75%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '12.47'
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  -10.98%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.5583'
$ws.Range('E47').NumberFormat = "@"
$ws.Range('E47').Value = '  -9.76%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '122.17'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  -6.45%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '1.940'
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  -9.44%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.06882'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  -6.02%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '1.175'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -4.22%  '
